{"js": "// Update the date line and the 25 multiplication problems to the new\n// values per the diff. Each original text is unique in the document, so a\n// simple exact-text search/replace (matchCase, no wildcards) is unambiguous.\nconst replacements = [\n  [\"2023-09-04 Monday\", \"2023-09-05 Tuesday\"],\n  [\"86\u00d749=\", \"36\u00d743=\"],\n  [\"17\u00d777=\", \"67\u00d711=\"],\n  [\"16\u00d748=\", \"16\u00d729=\"],\n  [\"60\u00d732=\", \"71\u00d763=\"],\n  [\"46\u00d779=\", \"38\u00d777=\"],\n  [\"81\u00d773=\", \"58\u00d796=\"],\n  [\"63\u00d788=\", \"71\u00d751=\"],\n  [\"76\u00d795=\", \"27\u00d793=\"],\n  [\"57\u00d746=\", \"46\u00d752=\"],\n  [\"48\u00d723=\", \"81\u00d764=\"],\n  [\"20\u00d769=\", \"62\u00d778=\"],\n  [\"21\u00d716=\", \"44\u00d760=\"],\n  [\"53\u00d713=\", \"45\u00d720=\"],\n  [\"57\u00d770=\", \"59\u00d782=\"],\n  [\"67\u00d753=\", \"49\u00d714=\"],\n  [\"63\u00d736=\", \"13\u00d772=\"],\n  [\"74\u00d784=\", \"81\u00d773=\"],\n  [\"98\u00d757=\", \"85\u00d726=\"],\n  [\"72\u00d736=\", \"50\u00d719=\"],\n  [\"11\u00d796=\", \"19\u00d756=\"],\n  [\"49\u00d727=\", \"95\u00d737=\"],\n  [\"79\u00d719=\", \"42\u00d783=\"],\n  [\"14\u00d728=\", \"86\u00d724=\"],\n  [\"25\u00d779=\", \"73\u00d757=\"],\n  [\"99\u00d722=\", \"94\u00d757=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 multiplication problems to the new\n# values per the diff. Each original text is unique in the document, so a\n# simple Find/Replace (MatchCase, no wildcards) over the whole story is\n# unambiguous.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-04 Monday\", \"2023-09-05 Tuesday\"),\n    @(\"86\u00d749=\", \"36\u00d743=\"),\n    @(\"17\u00d777=\", \"67\u00d711=\"),\n    @(\"16\u00d748=\", \"16\u00d729=\"),\n    @(\"60\u00d732=\", \"71\u00d763=\"),\n    @(\"46\u00d779=\", \"38\u00d777=\"),\n    @(\"81\u00d773=\", \"58\u00d796=\"),\n    @(\"63\u00d788=\", \"71\u00d751=\"),\n    @(\"76\u00d795=\", \"27\u00d793=\"),\n    @(\"57\u00d746=\", \"46\u00d752=\"),\n    @(\"48\u00d723=\", \"81\u00d764=\"),\n    @(\"20\u00d769=\", \"62\u00d778=\"),\n    @(\"21\u00d716=\", \"44\u00d760=\"),\n    @(\"53\u00d713=\", \"45\u00d720=\"),\n    @(\"57\u00d770=\", \"59\u00d782=\"),\n    @(\"67\u00d753=\", \"49\u00d714=\"),\n    @(\"63\u00d736=\", \"13\u00d772=\"),\n    @(\"74\u00d784=\", \"81\u00d773=\"),\n    @(\"98\u00d757=\", \"85\u00d726=\"),\n    @(\"72\u00d736=\", \"50\u00d719=\"),\n    @(\"11\u00d796=\", \"19\u00d756=\"),\n    @(\"49\u00d727=\", \"95\u00d737=\"),\n    @(\"79\u00d719=\", \"42\u00d783=\"),\n    @(\"14\u00d728=\", \"86\u00d724=\"),\n    @(\"25\u00d779=\", \"73\u00d757=\"),\n    @(\"99\u00d722=\", \"94\u00d757=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $oldText,\n        $true,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        1,\n        $false,\n        $newText,\n        2\n    ) | Out-Null\n}\n"}
